$d = $word.ActiveDocument

# Harvard paragraph: PursuingHappy entry is missing the "R:" label before "3/1)"
$d.Content.Find.Execute("I: 12/x, 3/1)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "I: 12/x, R:3/1)", 2)

# University of Kentucky paragraph: add "random_user " applicant before the dates
$d.Content.Find.Execute("University of Kentucky: (I: 11/12, R: 4/12)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "University of Kentucky: random_user (I: 11/12, R: 4/12)", 2)
